$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 updates - B2 and D2 are cleared (deleted), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.92012669354727072
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -0.3863929271143855

# Row 3 updates
$ws.Range("B3").Value = -1.3115821127436811
$ws.Range("C3").Value = 0.67302120142693267
$ws.Range("D3").Value = -0.61150265465160758
$ws.Range("E3").Value = 2.1505816834748517

# Update the selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
